# Generate Report for Handoff
#
# The "b.md" file has now had its own handoff package generated (rather
# than piggy-backing on "a.md" as a content duplicate), so the status
# rolls forward to "Ready for handoff" for both locales, the new handoff
# file name / timestamp for "b.md" is recorded, the row is flagged as no
# longer being a content duplicate, and a warning is attached noting the
# handback file that shipped is behind the newly generated source.

$wb = $excel.ActiveWorkbook

$statusReadyForHandoff = "Ready for handoff"
$overviewDateTime      = "2016-08-17 00:35:18"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3de9cb4f075a77cccc1155671ef7be9b5bb207ff/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aaa0abd2a4cd60641ec55d2a40a2512b5a742bb7/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is "b.md"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(3, 5).Value = $statusReadyForHandoff  # E3 - zh-cn status
$wsOverview.Cells.Item(3, 6).Value = $statusReadyForHandoff  # F3 - de-de status
$wsOverview.Cells.Item(3, 7).Value = $overviewDateTime       # G3 - Latest HO Xliff Generate Date

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is "b.md"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(3, 3).Value  = $statusReadyForHandoff                                         # C3 - Status
$wsZhCn.Cells.Item(3, 6).Value  = "False"                                                        # F3 - Content Duplicate
$wsZhCn.Cells.Item(3, 7).Value  = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"          # G3 - Latest Handoff File
$wsZhCn.Cells.Item(3, 8).Value  = "2016-08-17 00:35:13"                                          # H3 - Latest Handoff Datetime
$wsZhCn.Cells.Item(3, 16).Value = $errorDetail                                                    # P3 - Error Detail
$wsZhCn.Range("P:P").ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet: row 3 is "b.md"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(3, 3).Value  = $statusReadyForHandoff                                          # C3 - Status
$wsDeDe.Cells.Item(3, 6).Value  = "False"                                                         # F3 - Content Duplicate
$wsDeDe.Cells.Item(3, 7).Value  = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"           # G3 - Latest Handoff File
$wsDeDe.Cells.Item(3, 8).Value  = $overviewDateTime                                               # H3 - Latest Handoff Datetime
$wsDeDe.Cells.Item(3, 16).Value = $errorDetail                                                     # P3 - Error Detail
$wsDeDe.Range("P:P").ColumnWidth = 39.15

Write-Output "Applied handoff-report updates for b.md (zh-cn / de-de)."
